{"js": "// The document body contains one table of 100 arithmetic prompts such as\n// \"42-28=\" / \"1+39=\", laid out row-major (20 rows x 5 columns). The\n// commit swaps each prompt's text for a new prompt, cell by cell, in\n// table reading order (top-left to bottom-right). A few old prompts are\n// repeated (e.g. \"74-48=\" shows up twice) but map to *different* new\n// values depending on position, so cells are updated strictly in\n// document order rather than via a single global find/replace.\nconst REPLACEMENTS = [\n  [\"42-28=\", \"85-76=\"],\n  [\"42-33=\", \"53-3=\"],\n  [\"90-52=\", \"81-72=\"],\n  [\"98-24=\", \"17+5=\"],\n  [\"1+39=\", \"1+49=\"],\n  [\"94-57=\", \"69+15=\"],\n  [\"76-22=\", \"97-44=\"],\n  [\"52-3=\", \"73-67=\"],\n  [\"37-6=\", \"83-80=\"],\n  [\"8+53=\", \"32-23=\"],\n  [\"94-65=\", \"22+56=\"],\n  [\"75-39=\", \"34+19=\"],\n  [\"94-21=\", \"38+13=\"],\n  [\"44+2=\", \"13+74=\"],\n  [\"28+57=\", \"52+20=\"],\n  [\"99-77=\", \"89-30=\"],\n  [\"70+2=\", \"21+20=\"],\n  [\"80-22=\", \"87-14=\"],\n  [\"95-0=\", \"99-47=\"],\n  [\"41-7=\", \"95-64=\"],\n  [\"20+35=\", \"39+32=\"],\n  [\"70-42=\", \"53-43=\"],\n  [\"74-48=\", \"36+40=\"],\n  [\"72-3=\", \"76-26=\"],\n  [\"64-7=\", \"31+68=\"],\n  [\"21+22=\", \"60-18=\"],\n  [\"7+29=\", \"44-16=\"],\n  [\"32-5=\", \"16+32=\"],\n  [\"38+3=\", \"72+0=\"],\n  [\"97-68=\", \"93-0=\"],\n  [\"95-19=\", \"93-54=\"],\n  [\"53+14=\", \"80-32=\"],\n  [\"67-31=\", \"40-18=\"],\n  [\"64-0=\", \"54+27=\"],\n  [\"0+56=\", \"55+38=\"],\n  [\"1+30=\", \"23+9=\"],\n  [\"44+51=\", \"32-17=\"],\n  [\"7+14=\", \"88-86=\"],\n  [\"93+2=\", \"72-62=\"],\n  [\"86-75=\", \"16+10=\"],\n  [\"7+24=\", \"60-3=\"],\n  [\"59-19=\", \"52+45=\"],\n  [\"51-20=\", \"40+59=\"],\n  [\"96-68=\", \"72-70=\"],\n  [\"42-30=\", \"42+53=\"],\n  [\"10+35=\", \"98-63=\"],\n  [\"28-19=\", \"71-47=\"],\n  [\"36-24=\", \"75+11=\"],\n  [\"99-22=\", \"68+0=\"],\n  [\"9+39=\", \"15+15=\"],\n  [\"50-14=\", \"95-41=\"],\n  [\"4+61=\", \"54+12=\"],\n  [\"75-59=\", \"91-6=\"],\n  [\"66-46=\", \"4+52=\"],\n  [\"46+35=\", \"95-28=\"],\n  [\"75-42=\", \"24+59=\"],\n  [\"22+3=\", \"80+8=\"],\n  [\"46+44=\", \"91-79=\"],\n  [\"96-39=\", \"68+19=\"],\n  [\"75-19=\", \"63-21=\"],\n  [\"27-0=\", \"65-20=\"],\n  [\"69-66=\", \"95-86=\"],\n  [\"55-50=\", \"57+6=\"],\n  [\"76-71=\", \"90-43=\"],\n  [\"25-24=\", \"75-37=\"],\n  [\"15+22=\", \"1+87=\"],\n  [\"46+15=\", \"40-31=\"],\n  [\"93-64=\", \"32+37=\"],\n  [\"65-43=\", \"14+62=\"],\n  [\"41+6=\", \"55-33=\"],\n  [\"78-28=\", \"55+38=\"],\n  [\"29-25=\", \"26-9=\"],\n  [\"77-71=\", \"83-10=\"],\n  [\"98-52=\", \"68+0=\"],\n  [\"74-48=\", \"45+26=\"],\n  [\"33-2=\", \"50-2=\"],\n  [\"70-53=\", \"26-21=\"],\n  [\"86-22=\", \"62+14=\"],\n  [\"90-6=\", \"42+10=\"],\n  [\"83-23=\", \"10+7=\"],\n  [\"89-87=\", \"46-21=\"],\n  [\"7+89=\", \"8+76=\"],\n  [\"87-55=\", \"21+61=\"],\n  [\"49-4=\", \"64-35=\"],\n  [\"80+14=\", \"93-47=\"],\n  [\"93-79=\", \"61-50=\"],\n  [\"12+43=\", \"82-77=\"],\n  [\"35+28=\", \"97-54=\"],\n  [\"38+1=\", \"16-16=\"],\n  [\"85-74=\", \"15+79=\"],\n  [\"99-65=\", \"62-8=\"],\n  [\"35+35=\", \"37+26=\"],\n  [\"71-17=\", \"12+14=\"],\n  [\"18-16=\", \"73-66=\"],\n  [\"48+29=\", \"85-8=\"],\n  [\"53+24=\", \"97-61=\"],\n  [\"77-11=\", \"81-38=\"],\n  [\"4+87=\", \"88+3=\"],\n  [\"31+25=\", \"87-86=\"],\n  [\"43+33=\", \"20+40=\"],\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nif (tables.items.length === 0) {\n  throw new Error(\"Expected at least one table in the document body.\");\n}\n\nconst table = tables.items[0];\ntable.load(\"values\");\nawait context.sync();\n\nconst grid = table.values;\nconst flatCells = [];\nfor (let r = 0; r < grid.length; r++) {\n  for (let c = 0; c < grid[r].length; c++) {\n    flatCells.push({ r, c, text: grid[r][c] });\n  }\n}\n\nif (flatCells.length !== REPLACEMENTS.length) {\n  throw new Error(\n    `Expected ${REPLACEMENTS.length} table cells, found ${flatCells.length}.`\n  );\n}\n\nfor (let i = 0; i < REPLACEMENTS.length; i++) {\n  const [oldText, newText] = REPLACEMENTS[i];\n  const cellInfo = flatCells[i];\n  if (cellInfo.text !== oldText) {\n    throw new Error(\n      `Cell ${i} (row ${cellInfo.r}, col ${cellInfo.c}) expected \"${oldText}\" but found \"${cellInfo.text}\".`\n    );\n  }\n  table.getCell(cellInfo.r, cellInfo.c).value = newText;\n}\n\nawait context.sync();\n", "ps1": "# The document body contains one table of 100 arithmetic prompts such as\n# \"42-28=\" / \"1+39=\", laid out row-major (20 rows x 5 columns). The commit\n# swaps each prompt's text for a new prompt, cell by cell, in table\n# reading order (top-left to bottom-right). A few old prompts repeat\n# (e.g. \"74-48=\" shows up twice) but map to *different* new values\n# depending on position, so cells are updated strictly in document order\n# rather than via a single global find/replace.\n$replacements = @(\n    @(\"42-28=\", \"85-76=\"),\n    @(\"42-33=\", \"53-3=\"),\n    @(\"90-52=\", \"81-72=\"),\n    @(\"98-24=\", \"17+5=\"),\n    @(\"1+39=\", \"1+49=\"),\n    @(\"94-57=\", \"69+15=\"),\n    @(\"76-22=\", \"97-44=\"),\n    @(\"52-3=\", \"73-67=\"),\n    @(\"37-6=\", \"83-80=\"),\n    @(\"8+53=\", \"32-23=\"),\n    @(\"94-65=\", \"22+56=\"),\n    @(\"75-39=\", \"34+19=\"),\n    @(\"94-21=\", \"38+13=\"),\n    @(\"44+2=\", \"13+74=\"),\n    @(\"28+57=\", \"52+20=\"),\n    @(\"99-77=\", \"89-30=\"),\n    @(\"70+2=\", \"21+20=\"),\n    @(\"80-22=\", \"87-14=\"),\n    @(\"95-0=\", \"99-47=\"),\n    @(\"41-7=\", \"95-64=\"),\n    @(\"20+35=\", \"39+32=\"),\n    @(\"70-42=\", \"53-43=\"),\n    @(\"74-48=\", \"36+40=\"),\n    @(\"72-3=\", \"76-26=\"),\n    @(\"64-7=\", \"31+68=\"),\n    @(\"21+22=\", \"60-18=\"),\n    @(\"7+29=\", \"44-16=\"),\n    @(\"32-5=\", \"16+32=\"),\n    @(\"38+3=\", \"72+0=\"),\n    @(\"97-68=\", \"93-0=\"),\n    @(\"95-19=\", \"93-54=\"),\n    @(\"53+14=\", \"80-32=\"),\n    @(\"67-31=\", \"40-18=\"),\n    @(\"64-0=\", \"54+27=\"),\n    @(\"0+56=\", \"55+38=\"),\n    @(\"1+30=\", \"23+9=\"),\n    @(\"44+51=\", \"32-17=\"),\n    @(\"7+14=\", \"88-86=\"),\n    @(\"93+2=\", \"72-62=\"),\n    @(\"86-75=\", \"16+10=\"),\n    @(\"7+24=\", \"60-3=\"),\n    @(\"59-19=\", \"52+45=\"),\n    @(\"51-20=\", \"40+59=\"),\n    @(\"96-68=\", \"72-70=\"),\n    @(\"42-30=\", \"42+53=\"),\n    @(\"10+35=\", \"98-63=\"),\n    @(\"28-19=\", \"71-47=\"),\n    @(\"36-24=\", \"75+11=\"),\n    @(\"99-22=\", \"68+0=\"),\n    @(\"9+39=\", \"15+15=\"),\n    @(\"50-14=\", \"95-41=\"),\n    @(\"4+61=\", \"54+12=\"),\n    @(\"75-59=\", \"91-6=\"),\n    @(\"66-46=\", \"4+52=\"),\n    @(\"46+35=\", \"95-28=\"),\n    @(\"75-42=\", \"24+59=\"),\n    @(\"22+3=\", \"80+8=\"),\n    @(\"46+44=\", \"91-79=\"),\n    @(\"96-39=\", \"68+19=\"),\n    @(\"75-19=\", \"63-21=\"),\n    @(\"27-0=\", \"65-20=\"),\n    @(\"69-66=\", \"95-86=\"),\n    @(\"55-50=\", \"57+6=\"),\n    @(\"76-71=\", \"90-43=\"),\n    @(\"25-24=\", \"75-37=\"),\n    @(\"15+22=\", \"1+87=\"),\n    @(\"46+15=\", \"40-31=\"),\n    @(\"93-64=\", \"32+37=\"),\n    @(\"65-43=\", \"14+62=\"),\n    @(\"41+6=\", \"55-33=\"),\n    @(\"78-28=\", \"55+38=\"),\n    @(\"29-25=\", \"26-9=\"),\n    @(\"77-71=\", \"83-10=\"),\n    @(\"98-52=\", \"68+0=\"),\n    @(\"74-48=\", \"45+26=\"),\n    @(\"33-2=\", \"50-2=\"),\n    @(\"70-53=\", \"26-21=\"),\n    @(\"86-22=\", \"62+14=\"),\n    @(\"90-6=\", \"42+10=\"),\n    @(\"83-23=\", \"10+7=\"),\n    @(\"89-87=\", \"46-21=\"),\n    @(\"7+89=\", \"8+76=\"),\n    @(\"87-55=\", \"21+61=\"),\n    @(\"49-4=\", \"64-35=\"),\n    @(\"80+14=\", \"93-47=\"),\n    @(\"93-79=\", \"61-50=\"),\n    @(\"12+43=\", \"82-77=\"),\n    @(\"35+28=\", \"97-54=\"),\n    @(\"38+1=\", \"16-16=\"),\n    @(\"85-74=\", \"15+79=\"),\n    @(\"99-65=\", \"62-8=\"),\n    @(\"35+35=\", \"37+26=\"),\n    @(\"71-17=\", \"12+14=\"),\n    @(\"18-16=\", \"73-66=\"),\n    @(\"48+29=\", \"85-8=\"),\n    @(\"53+24=\", \"97-61=\"),\n    @(\"77-11=\", \"81-38=\"),\n    @(\"4+87=\", \"88+3=\"),\n    @(\"31+25=\", \"87-86=\"),\n    @(\"43+33=\", \"20+40=\")\n)\n\n$d = $word.ActiveDocument\nif ($d.Tables.Count -lt 1) {\n    throw \"Expected at least one table in the document.\"\n}\n$t = $d.Tables.Item(1)\n\n$rows = $t.Rows.Count\n$cols = $t.Columns.Count\nif (($rows * $cols) -ne $replacements.Length) {\n    throw \"Expected $($replacements.Length) table cells, found $($rows * $cols).\"\n}\n\n$i = 0\nfor ($r = 1; $r -le $rows; $r++) {\n    for ($c = 1; $c -le $cols; $c++) {\n        $pair = $replacements[$i]\n        $oldText = $pair[0]\n        $newText = $pair[1]\n        $cell = $t.Cell($r, $c)\n        $cellRange = $cell.Range\n        $current = $cellRange.Text.TrimEnd([char]13, [char]7)\n        if ($current -ne $oldText) {\n            throw \"Cell ($r,$c) expected '$oldText' but found '$current'.\"\n        }\n        $cellRange.Text = $newText\n        $i++\n    }\n}\n"}
